# Apply the signoff-sheet edit described by the commit:
# "add additional signoffs, update flow chart, change some new courses to
#  revised courses"
#
# Concretely this collapses 5 checkbox content controls down to a single
# run each (dropping a stray trailing empty run left over in their
# sdtContent), merges the "Rationale: " run with the course-change text
# that follows it, and turns the "Potential Impact on Other Departments"
# answer from "NONE" into an underlined "Computer Electives List".

$d = $word.ActiveDocument

# --- 1) Collapse the checkbox content controls --------------------------
# Re-assigning ContentControl.Range.Text replaces *all* the runs inside
# the sdtContent with a single run (using the formatting of the content
# control's first run) carrying just the given text - which is exactly
# what's needed to drop the extra empty trailing <w:r> the diff removes.
$checkboxIds = @(
    @{ Id = -647354619;  Glyph = [char]0x2610 },  # unchecked
    @{ Id = -1850945197; Glyph = [char]0x2612 },  # checked
    @{ Id = -570272621;  Glyph = [char]0x2612 },  # checked
    @{ Id = -567805947;  Glyph = [char]0x2612 },  # checked
    @{ Id = 1345512495;  Glyph = [char]0x2612 }   # checked
)

for ($i = 1; $i -le $d.ContentControls.Count; $i++) {
    $cc = $d.ContentControls.Item($i)
    foreach ($entry in $checkboxIds) {
        if ($cc.ID -eq $entry.Id) {
            $cc.Range.Text = $entry.Glyph
        }
    }
}

# --- 2) Merge "Rationale: " with the course-change text ------------------
# Both runs already share the same formatting, so re-running Find/Replace
# over them merges the two runs into one (matching the diff).
$rng = $d.Content
$rng.Find.Execute("Rationale: ", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "Rationale: ", 2) | Out-Null

# --- 3) "NONE" -> "Computer Electives List" with underline ---------------
# Only the first "Potential Impact on Other Departments" answer changes;
# scope the Find to just that part of the document so the later
# "Potential Resource/Financial Requirements" NONE is left untouched.
$rng2 = $d.Content
$found = $rng2.Find.Execute("Potential Impact on Other Departments: NONE")
if ($found) {
    $target = $d.Range($rng2.Start, $rng2.End)
    $target.Find.Execute("NONE", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "Computer Electives List", 2) | Out-Null
}

$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Computer Electives List")
if ($found3) {
    $rng3.Font.Underline = 1
}
